# Apply header renames and title-case fixes to municipality/state names
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
$ws.Range("B20").Value = "Amatenango De La Frontera"
$ws.Range("B24").Value = "Benemérito De Las Américas"
$ws.Range("B28").Value = "Chiapa De Corzo"
$ws.Range("B30").Value = "Comitán De Domínguez"
$ws.Range("B50").Value = "San Cristóbal De Las Casas"
$ws.Range("B80").Value = "Hidalgo Del Parral"
$ws.Range("B83").Value = "San Francisco De Borja"
$ws.Range("B95").Value = "San Juan De Sabinas"
$ws.Range("A104").Value = "Ciudad De México"
$ws.Range("B108").Value = "Cuajimalpa De Morelos"
$ws.Range("B126").Value = "Pánuco De Coronado"
$ws.Range("B127").Value = "San Juan De Guadalupe"
$ws.Range("A132").Value = "Estado De México"
$ws.Range("B132").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B134").Value = "Atizapán De Zaragoza"
$ws.Range("B140").Value = "Ecatepec De Morelos"
$ws.Range("B142").Value = "Ixtapan De La Sal"
$ws.Range("B148").Value = "Naucalpan De Juárez"
$ws.Range("B150").Value = "San Felipe Del Progreso"
$ws.Range("B151").Value = "San Martín De Las Pirámides"
$ws.Range("B158").Value = "Tlalnepantla De Baz"
$ws.Range("B167").Value = "San Miguel De Allende"
$ws.Range("B168").Value = "Apaseo El Alto"
$ws.Range("B169").Value = "Apaseo El Grande"
$ws.Range("B173").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B177").Value = "Jaral Del Progreso"
$ws.Range("B187").Value = "San Francisco Del Rincón"
$ws.Range("B189").Value = "San Luis De La Paz"
$ws.Range("B190").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B191").Value = "Silao De La Victoria"
$ws.Range("B194").Value = "Valle De Santiago"
$ws.Range("B199").Value = "Acapulco De Juárez"
$ws.Range("B201").Value = "Ajuchitlán Del Progreso"
$ws.Range("B204").Value = "Atlamajalcingo Del Monte"
$ws.Range("B206").Value = "Atoyac De Álvarez"
$ws.Range("B207").Value = "Ayutla De Los Libres"
$ws.Range("B209").Value = "Chilapa De Álvarez"
$ws.Range("B210").Value = "Chilpancingo De Los Bravo"
$ws.Range("B211").Value = "Coyuca De Benítez"
$ws.Range("B212").Value = "Cutzamala De Pinzón"
$ws.Range("B213").Value = "Huitzuco De Los Figueroa"
$ws.Range("B214").Value = "Iguala De La Independencia"
$ws.Range("B216").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B219").Value = "Mártir De Cuilapan"
$ws.Range("B226").Value = "Taxco De Alarcón"
$ws.Range("B227").Value = "Técpan De Galeana"
$ws.Range("B228").Value = "Tepecoacuilco De Trujano"
$ws.Range("B232").Value = "Tlapa De Comonfort"
$ws.Range("B249").Value = "Pachuca De Soto"
$ws.Range("B251").Value = "Progreso De Obregón"
$ws.Range("B254").Value = "Tenango De Doria"
$ws.Range("B255").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B259").Value = "Tula De Allende"
$ws.Range("B263").Value = "Autlán De Navarro"
$ws.Range("B268").Value = "Encarnación De Díaz"
$ws.Range("B272").Value = "La Manzanilla De La Paz"
$ws.Range("B273").Value = "Lagos De Moreno"
$ws.Range("B278").Value = "San Juan De Los Lagos"
$ws.Range("B279").Value = "San Miguel El Alto"
$ws.Range("B328").Value = "Tlaltizapán De Zapata"
$ws.Range("B333").Value = "Ixtlán Del Río"
$ws.Range("B345").Value = "Mier Y Noriega"
$ws.Range("B348").Value = "San Nicolás De Los Garza"
$ws.Range("B351").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B353").Value = "Coicoyán De Las Flores"
$ws.Range("B355").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B356").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B357").Value = "Ixtlán De Juárez"
$ws.Range("B358").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B360").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B367").Value = "San José Del Progreso"
$ws.Range("B396").Value = "Santo Domingo De Morelos"
$ws.Range("B401").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B402").Value = "Tataltepec De Valdés"
$ws.Range("B403").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B404").Value = "Tlacolula De Matamoros"
$ws.Range("B405").Value = "Villa De Tututepec"
$ws.Range("B415").Value = "Izúcar De Matamoros"
$ws.Range("B422").Value = "San Salvador El Seco"
$ws.Range("B424").Value = "Tecali De Herrera"
$ws.Range("B428").Value = "Tetela De Ocampo"
$ws.Range("B442").Value = "Amealco De Bonfil"
$ws.Range("B444").Value = "Cadereyta De Montes"
$ws.Range("B447").Value = "Jalpan De Serra"
$ws.Range("B448").Value = "Landa De Matamoros"
$ws.Range("B450").Value = "Pinal De Amoles"
$ws.Range("B453").Value = "San Juan Del Río"
$ws.Range("B462").Value = "Cerro De San Pedro"
$ws.Range("B463").Value = "Ciudad Del Maíz"
$ws.Range("B475").Value = "Santa María Del Río"
$ws.Range("B482").Value = "Villa De Guadalupe"
$ws.Range("B483").Value = "Villa De Ramos"
$ws.Range("B525").Value = "Soto La Marina"
$ws.Range("B539").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B541").Value = "Amatlán De Los Reyes"
$ws.Range("B545").Value = "Boca Del Río"
$ws.Range("B547").Value = "Castillo De Teayo"
$ws.Range("B556").Value = "Cosamaloapan De Carpio"
$ws.Range("B557").Value = "Cosautlán De Carvajal"
$ws.Range("B567").Value = "Hueyapan De Ocampo"
$ws.Range("B568").Value = "Ignacio De La Llave"
$ws.Range("B576").Value = "Lerdo De Tejada"
$ws.Range("B577").Value = "Martínez De La Torre"
$ws.Range("B581").Value = "Mixtla De Altamirano"
$ws.Range("B586").Value = "Paso Del Macho"
$ws.Range("B587").Value = "Poza Rica De Hidalgo"
$ws.Range("B594").Value = "Sayula De Alemán"
$ws.Range("B595").Value = "Soledad De Doblado"
$ws.Range("B609").Value = "Vega De Alatorre"
$ws.Range("B626").Value = "Nochistlán De Mejía"
$ws.Range("B632").Value = "Tlaltenango De Sánchez Román"

# Minor floating point precision corrections (re-save artifact)
$ws.Range("D18").Value = 0.009350649350649352
$ws.Range("D519").Value = 0.009350649350649352
$ws.Range("D67").Value = 0.09402597402597404

# Remove trailing footnote/metadata rows (637-641); row 636 stays blank/absent
$ws.Range("A637:A641").EntireRow.Delete()
